$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The Cypher query text that now lives in cell A2 (added as a new shared string).
$query = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.ethnicity IN ['UNKNOWN'] RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"

# Populate A2 with the query (cell previously was blank but carried the wrap-text style).
$ws.Range("A2").Value = $query

# Give the row enough height to show the wrapped, multi-line query text.
$ws.Rows(2).RowHeight = 87

# Move the selection/active cell to A2 and scroll the view back to the top-left.
$ws.Range("A2").Select()
